# Auto-generated edit script: updates cached market-data cells (columns H-N)
# across multiple sheets to match the scheduled runner refresh described in the diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (54 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    @{Addr="H28"; Val=845.13336},
    @{Addr="I28"; Val=656.4545000000001},
    @{Addr="K28"; Val=656.4545000000001},
    @{Addr="M28"; Val=-171.4545000000001},
    @{Addr="H62"; Val=6236.1934},
    @{Addr="I62"; Val=5740.68},
    @{Addr="J62"; Val=8300.833000000001},
    @{Addr="K62"; Val=5740.68},
    @{Addr="L62"; Val=8300.833000000001},
    @{Addr="M62"; Val=-5116.68},
    @{Addr="N62"; Val=-9548.833000000001},
    @{Addr="H65"; Val=6236.1934},
    @{Addr="I65"; Val=5740.68},
    @{Addr="J65"; Val=8300.833000000001},
    @{Addr="K65"; Val=28703.4},
    @{Addr="L65"; Val=41504.165},
    @{Addr="M65"; Val=-25583.4},
    @{Addr="N65"; Val=-47744.165},
    @{Addr="H86"; Val=2864.3333},
    @{Addr="I86"; Val=2787.2},
    @{Addr="J86"; Val=3250},
    @{Addr="K86"; Val=2787.2},
    @{Addr="L86"; Val=3250},
    @{Addr="M86"; Val=-1664.2},
    @{Addr="N86"; Val=-5496},
    @{Addr="H89"; Val=2864.3333},
    @{Addr="I89"; Val=2787.2},
    @{Addr="J89"; Val=3250},
    @{Addr="K89"; Val=13936},
    @{Addr="L89"; Val=16250},
    @{Addr="M89"; Val=-8320},
    @{Addr="N89"; Val=-27482},
    @{Addr="H98"; Val=1880.579},
    @{Addr="I98"; Val=1602.5294},
    @{Addr="J98"; Val=4244},
    @{Addr="K98"; Val=1602.5294},
    @{Addr="L98"; Val=4244},
    @{Addr="M98"; Val=-104.5293999999999},
    @{Addr="N98"; Val=-7240},
    @{Addr="H106"; Val=3722.5454},
    @{Addr="I106"; Val=3661},
    @{Addr="K106"; Val=3661},
    @{Addr="M106"; Val=-3030},
    @{Addr="H122"; Val=1880.579},
    @{Addr="I122"; Val=1602.5294},
    @{Addr="J122"; Val=4244},
    @{Addr="K122"; Val=4807.5882},
    @{Addr="L122"; Val=12732},
    @{Addr="M122"; Val=-2357.5882},
    @{Addr="N122"; Val=-17632},
    @{Addr="H141"; Val=3074.5},
    @{Addr="I141"; Val=3074.5},
    @{Addr="K141"; Val=9223.5},
    @{Addr="M141"; Val=-4043.5}
)
foreach ($u in $updates) {
    $ws.Range($u.Addr).Value = $u.Val
}

# --- Sheet: ARM (30 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    @{Addr="H32"; Val=17419.092},
    @{Addr="I32"; Val=17822.812},
    @{Addr="K32"; Val=17822.812},
    @{Addr="M32"; Val=-17535.812},
    @{Addr="H45"; Val=3916.2354},
    @{Addr="I45"; Val=2262.4546},
    @{Addr="J45"; Val=6948.1665},
    @{Addr="K45"; Val=2262.4546},
    @{Addr="L45"; Val=6948.1665},
    @{Addr="M45"; Val=-1885.4546},
    @{Addr="N45"; Val=-7702.1665},
    @{Addr="H74"; Val=226858.62},
    @{Addr="I74"; Val=261769.12},
    @{Addr="K74"; Val=261769.12},
    @{Addr="M74"; Val=-260895.12},
    @{Addr="H77"; Val=226858.62},
    @{Addr="I77"; Val=261769.12},
    @{Addr="K77"; Val=1308845.6},
    @{Addr="M77"; Val=-1304477.6},
    @{Addr="H102"; Val=6116.6665},
    @{Addr="I102"; Val=6637.5},
    @{Addr="J102"; Val=1950},
    @{Addr="K102"; Val=6637.5},
    @{Addr="L102"; Val=1950},
    @{Addr="M102"; Val=-5015.5},
    @{Addr="N102"; Val=-5194},
    @{Addr="H132"; Val=2305.5293},
    @{Addr="J132"; Val=2171.875},
    @{Addr="L132"; Val=6515.625},
    @{Addr="N132"; Val=-11575.625}
)
foreach ($u in $updates) {
    $ws.Range($u.Addr).Value = $u.Val
}

# --- Sheet: BSM (19 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    @{Addr="H86"; Val=1707.7142},
    @{Addr="I86"; Val=1716.4375},
    @{Addr="K86"; Val=1716.4375},
    @{Addr="M86"; Val=-593.4375},
    @{Addr="H89"; Val=1707.7142},
    @{Addr="I89"; Val=1716.4375},
    @{Addr="K89"; Val=8582.1875},
    @{Addr="M89"; Val=-2966.1875},
    @{Addr="H99"; Val=1663.375},
    @{Addr="I99"; Val=1543.8572},
    @{Addr="K99"; Val=1543.8572},
    @{Addr="M99"; Val=-45.85719999999992},
    @{Addr="H134"; Val=2289.195},
    @{Addr="I134"; Val=1828.4839},
    @{Addr="J134"; Val=3717.4},
    @{Addr="K134"; Val=5485.4517},
    @{Addr="L134"; Val=11152.2},
    @{Addr="M134"; Val=-2950.4517},
    @{Addr="N134"; Val=-16222.2}
)
foreach ($u in $updates) {
    $ws.Range($u.Addr).Value = $u.Val
}

# --- Sheet: CRP (14 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    @{Addr="H31"; Val=3706706.2},
    @{Addr="I31"; Val=10001126},
    @{Addr="J31"; Val=4106.4707},
    @{Addr="K31"; Val=10001126},
    @{Addr="L31"; Val=4106.4707},
    @{Addr="M31"; Val=-10000831},
    @{Addr="N31"; Val=-4696.4707},
    @{Addr="H34"; Val=3706706.2},
    @{Addr="I34"; Val=10001126},
    @{Addr="J34"; Val=4106.4707},
    @{Addr="K34"; Val=10001126},
    @{Addr="L34"; Val=4106.4707},
    @{Addr="M34"; Val=-10000924},
    @{Addr="N34"; Val=-4510.4707}
)
foreach ($u in $updates) {
    $ws.Range($u.Addr).Value = $u.Val
}

# --- Sheet: CUL (29 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    @{Addr="H4"; Val=61333364},
    @{Addr="I4"; Val=45834656},
    @{Addr="J4"; Val=93522984},
    @{Addr="K4"; Val=137503968},
    @{Addr="L4"; Val=280568952},
    @{Addr="M4"; Val=-137503856},
    @{Addr="N4"; Val=-280569176},
    @{Addr="H8"; Val=157.36363},
    @{Addr="I8"; Val=157.36363},
    @{Addr="K8"; Val=472.09089},
    @{Addr="M8"; Val=-333.09089},
    @{Addr="H60"; Val=2461.125},
    @{Addr="I60"; Val=547.25},
    @{Addr="J60"; Val=4375},
    @{Addr="K60"; Val=1641.75},
    @{Addr="L60"; Val=13125},
    @{Addr="M60"; Val=-1390.75},
    @{Addr="N60"; Val=-13627},
    @{Addr="H92"; Val=1716.4},
    @{Addr="I92"; Val=1770.5},
    @{Addr="J92"; Val=1500},
    @{Addr="K92"; Val=5311.5},
    @{Addr="L92"; Val=4500},
    @{Addr="M92"; Val=-4063.5},
    @{Addr="N92"; Val=-6996},
    @{Addr="H107"; Val=1341.4706},
    @{Addr="I107"; Val=2202.3333},
    @{Addr="K107"; Val=6606.999899999999},
    @{Addr="M107"; Val=-4686.999899999999}
)
foreach ($u in $updates) {
    $ws.Range($u.Addr).Value = $u.Val
}

# --- Sheet: GSM (16 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    @{Addr="H58"; Val=37644.445},
    @{Addr="I58"; Val=18800},
    @{Addr="K58"; Val=18800},
    @{Addr="M58"; Val=-18523},
    @{Addr="H102"; Val=16186.404},
    @{Addr="I102"; Val=22549},
    @{Addr="K102"; Val=22549},
    @{Addr="M102"; Val=-20927},
    @{Addr="H122"; Val=3070.524},
    @{Addr="I122"; Val=2764.75},
    @{Addr="K122"; Val=8294.25},
    @{Addr="M122"; Val=-5844.25},
    @{Addr="H126"; Val=3272.0557},
    @{Addr="J126"; Val=4444},
    @{Addr="L126"; Val=13332},
    @{Addr="N126"; Val=-18272}
)
foreach ($u in $updates) {
    $ws.Range($u.Addr).Value = $u.Val
}

# --- Sheet: LTW (41 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    @{Addr="H7"; Val=3800.8},
    @{Addr="J7"; Val=7005},
    @{Addr="L7"; Val=7005},
    @{Addr="N7"; Val=-7229},
    @{Addr="H40"; Val=2282.75},
    @{Addr="I40"; Val=1778.3182},
    @{Addr="K40"; Val=1778.3182},
    @{Addr="M40"; Val=-1642.3182},
    @{Addr="H82"; Val=1176.0454},
    @{Addr="I82"; Val=1064.5454},
    @{Addr="K82"; Val=1064.5454},
    @{Addr="M82"; Val=-703.5454},
    @{Addr="H85"; Val=1176.0454},
    @{Addr="I85"; Val=1064.5454},
    @{Addr="K85"; Val=1064.5454},
    @{Addr="M85"; Val=183.4546},
    @{Addr="H88"; Val=54285.715},
    @{Addr="I88"; Val=25000},
    @{Addr="J88"; Val=59166.668},
    @{Addr="K88"; Val=25000},
    @{Addr="L88"; Val=59166.668},
    @{Addr="M88"; Val=-24572},
    @{Addr="N88"; Val=-60022.668},
    @{Addr="H91"; Val=54285.715},
    @{Addr="I91"; Val=25000},
    @{Addr="J91"; Val=59166.668},
    @{Addr="K91"; Val=25000},
    @{Addr="L91"; Val=59166.668},
    @{Addr="M91"; Val=-23518},
    @{Addr="N91"; Val=-62130.668},
    @{Addr="H126"; Val=3800.8},
    @{Addr="J126"; Val=7005},
    @{Addr="L126"; Val=21015},
    @{Addr="N126"; Val=-25955},
    @{Addr="H132"; Val=3453.5},
    @{Addr="I132"; Val=3205.5293},
    @{Addr="J132"; Val=4296.6},
    @{Addr="K132"; Val=9616.5879},
    @{Addr="L132"; Val=12889.8},
    @{Addr="M132"; Val=-7086.5879},
    @{Addr="N132"; Val=-17949.8}
)
foreach ($u in $updates) {
    $ws.Range($u.Addr).Value = $u.Val
}

# --- Sheet: WVR (29 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    @{Addr="H26"; Val=20499.666},
    @{Addr="I26"; Val=26999.5},
    @{Addr="K26"; Val=26999.5},
    @{Addr="M26"; Val=-26706.5},
    @{Addr="H81"; Val=17249.533},
    @{Addr="I81"; Val=21995},
    @{Addr="J81"; Val=4199.5},
    @{Addr="K81"; Val=43990},
    @{Addr="L81"; Val=8399},
    @{Addr="M81"; Val=-42929},
    @{Addr="N81"; Val=-10521},
    @{Addr="H84"; Val=17249.533},
    @{Addr="I84"; Val=21995},
    @{Addr="J84"; Val=4199.5},
    @{Addr="K84"; Val=219950},
    @{Addr="L84"; Val=41995},
    @{Addr="M84"; Val=-214646},
    @{Addr="N84"; Val=-52603},
    @{Addr="H122"; Val=61477.547},
    @{Addr="I122"; Val=69979},
    @{Addr="J122"; Val=7635},
    @{Addr="K122"; Val=209937},
    @{Addr="L122"; Val=22905},
    @{Addr="M122"; Val=-207487},
    @{Addr="N122"; Val=-27805},
    @{Addr="H126"; Val=281585.78},
    @{Addr="I126"; Val=2948.6667},
    @{Addr="K126"; Val=8846.000100000001},
    @{Addr="M126"; Val=-6376.000100000001}
)
foreach ($u in $updates) {
    $ws.Range($u.Addr).Value = $u.Val
}

Write-Output "Applied 232 cell updates across 8 sheets."